# Insert a new weekly price-report row (row 7) for "Hortaliza, Terminal
# Hortofrutícola Agro Chillán - Espinaca". Existing rows 7-12 shift down to
# become rows 8-13 (their data is unchanged), and the new row 7 is filled
# with the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data rows down by inserting a blank row at row 7.
$ws.Rows(7).Insert()

# Fill in the new row 7 with this week's record.
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = 44810
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 100112012
$ws.Range("G7").Value = "Espinaca"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 7500
$ws.Range("N7").Value = "`$/cuna 10 kilos"
$ws.Range("O7").Value = "Provincia de Diguillín"
$ws.Range("P7").Value = 750
$ws.Range("Q7").Value = 10
$ws.Range("R7").Value = "Hortaliza"
